$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4-6 (old entries for myemail3@, myemail5@, myemail@ with pwd 1)
$ws.Range("A4:D6").EntireRow.Delete() | Out-Null

# Update row 2: new verification entry (keep numeric/date-looking text as plain text,
# same as the Node.js backend's original inline-string writes)
$ws.Range("A2").Value = "myemail34@gmail.com"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "134"
$ws.Range("B2").ClearFormats() | Out-Null
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2024-11-17"
$ws.Range("C2").ClearFormats() | Out-Null
$ws.Range("D2").Value = "21:58:46"

# Update row 3: new verification entry
$ws.Range("A3").Value = "myemail@gmail.com"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "12"
$ws.Range("B3").ClearFormats() | Out-Null
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2024-11-17"
$ws.Range("C3").ClearFormats() | Out-Null
$ws.Range("D3").Value = "22:13:21"

# Widen column A to fit new content
$ws.Columns.Item(1).ColumnWidth = 22.140625

# Update selection to reflect last-used range
$ws.Range("A3:D3").Select() | Out-Null
